# Add profile configuration and bug fix
#
# 1) Workspace sheet: the "Business line" default value changes from
#    "Imprese Domestico" to "Corporate".
# 2) Selection / active-tab bookkeeping moves from the
#    "r Workspace_TargetVariable" sheet back to the main "Workspace" sheet.

$wb = $excel.ActiveWorkbook

$wsWorkspace = $wb.Worksheets.Item("Workspace")
$wsTargetVar = $wb.Worksheets.Item("r Workspace_TargetVariable")

# --- data fix: replace "Imprese Domestico" with "Corporate" -----------------
$wsWorkspace.Range("C3").Value = "Corporate"
$wsWorkspace.Range("D3").Value = "Corporate"

# --- view/selection bookkeeping ---------------------------------------------
# Move the selection on the (soon to be inactive) TargetVariable sheet first,
# so that activating/selecting the Workspace sheet afterwards is what sticks
# as the workbook's active tab.
$wsTargetVar.Activate()
$wsTargetVar.Range("D6").Select()

$wsWorkspace.Activate()
$wsWorkspace.Range("D3").Select()
